$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before the existing "Late" column (old column N),
# pushing "Late" to column O and "Outstanding" to column Q.
$ws.Columns("N").Insert() | Out-Null

# Match column N's width to its neighbouring "Waived" column (M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment Schedule" the active sheet (was "Transactions"), and
# update its selection.
$ws.Activate()
$ws.Range("T10").Select() | Out-Null
